$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "06/04/2025"
$ws.Cells.Item(20, 2).Value = 0.0004723999999999996
$ws.Cells.Item(20, 3).Value = 104784.081287045
$ws.Cells.Item(20, 4).Value = 49.5
